# "Calculo Mermelada" workbook update.
# Commit: "backend y frotend funcionando correctamente y pasando los
# reslutados correctos" — the real "Manzanas" price (I3, labeled "Precio")
# was corrected from 500 to 436. Every other cell that differs in the
# target file (N3, I4, N5, ...) is a formula that already depends on I3,
# so Excel's own recalculation engine produces the new values for us.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Correct the "Precio" input used by the right-hand "Teorico Manzanas"
# block (was 500, should be 436). The runtime auto-recalculates after the
# script runs, so every dependent formula (N3, I4, N5, ...) picks up the
# corrected value automatically before the workbook is saved.
$ws.Range("I3").Value = 436

# Leave the selection where the author left it when saving.
$ws.Range("I20").Select() | Out-Null
